$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.366.83'
$ws.Range('E2').Value = '  +4.80%  '
$ws.Range('D3').Value = '2.033.97'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.22'
$ws.Range('E5').Value = '  +5.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.644'
$ws.Range('E6').Value = '  +2.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.70'
$ws.Range('E7').Value = '  +16.45%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.375'
$ws.Range('E9').Value = '  +6.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.79'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0751'
$ws.Range('E11').Value = '  +4.23%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.905'
$ws.Range('E13').Value = '  +4.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.07'
$ws.Range('E14').Value = '  +9.32%  '
$ws.Range('D15').Value = '2.330.70'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.54'
$ws.Range('E16').Value = '  +7.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.27'
$ws.Range('E17').Value = '  +20.29%  '
$ws.Range('D18').Value = '2.027.49'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Value = '37.083.04'
$ws.Range('E19').Value = '  +4.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.09'
$ws.Range('E20').Value = '  +4.93%  '
$ws.Range('D21').Value = '0.0₃0871'
$ws.Range('E21').Value = '  +5.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.31'
$ws.Range('E22').Value = '  +6.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.84'
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.77'
$ws.Range('E24').Value = '  +25.42%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').Value = '  +6.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.43'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.79'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('E30').Value = '  +2.65%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.113'
$ws.Range('E31').Value = '  +29.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.22'
$ws.Range('E32').Value = '  +8.33%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.16'
$ws.Range('E33').Value = '  +9.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.67'
$ws.Range('E34').Value = '  +10.99%  '
$ws.Range('E35').Value = '  +5.55%  '
$ws.Range('E36').Value = '  +15.60%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.00'
$ws.Range('E38').Value = '  +25.05%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.80'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('E40').Value = '  +20.20%  '
$ws.Range('E41').Value = '  +5.18%  '
$ws.Range('E42').Value = '  +1.80%  '
$ws.Range('E43').Value = '  +25.81%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0218'
$ws.Range('E44').Value = '  +4.40%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.14'
$ws.Range('E45').Value = '  +6.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.08'
$ws.Range('E46').Value = '  +12.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.92'
$ws.Range('E47').Value = '  +12.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '94.92'
$ws.Range('E48').Value = '  +6.96%  '
$ws.Range('D49').Value = '1.428.83'
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.48'
$ws.Range('E51').Value = '  +6.74%  '
